# Insert a new data row at row 48 (pushing existing rows 48-63 down to 49-64)
# and populate it with the latest weekly price observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 48:63 down by one to make room for the new observation.
$ws.Rows("48:48").Insert()

# Populate the newly inserted row 48 with the new weekly record.
$ws.Range("A48").Value = 10
$ws.Range("B48").Value = "Vega Modelo de Temuco"
$ws.Range("C48").Value = "La Araucanía"
$ws.Range("D48").Value = 44704
$ws.Range("D48").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E48").Value = 9
$ws.Range("F48").Value = 300000001
$ws.Range("G48").Value = "Rabanito"
$ws.Range("H48").Value = "Sin especificar"
$ws.Range("I48").Value = "Primera"
$ws.Range("J48").Value = 40
$ws.Range("K48").Value = 7000
$ws.Range("L48").Value = 7000
$ws.Range("M48").Value = 7000
$ws.Range("N48").Value = "$/docena de paquetes"
$ws.Range("O48").Value = "Provincia de Cautín"
$ws.Range("P48").Value = 583
$ws.Range("Q48").Value = 12
$ws.Range("R48").Value = "Hortaliza"
